$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6892867684364319
$ws.Range("B1").Value = 2.170462369918823
$ws.Range("C1").Value = 5.053410530090332
$ws.Range("D1").Value = 3.018070220947266
$ws.Range("E1").Value = 0.7366823554039001
